$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" column (H) -------------------------------------------------
# Header text, matching the style already used by the other header cells
# (B1:G1 use style index 1 -> bold / bordered / centered).
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Block 1 (rows 2-11): Control 30/11/3/38/29 -> Label 0, MDD 41/8/15/16/33 -> Label 1
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# Block 2 (rows 12-21): same pattern
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1

# --- Refit results: updated Prediction / Error values for rows 4-11 --------
$ws.Range("D4").Value = 0.4330680414947439
$ws.Range("E4").Value = 0.4330680414947439

$ws.Range("D5").Value = 0.5990444533002196
$ws.Range("E5").Value = 0.5990444533002196

$ws.Range("D6").Value = 0.469287657961066
$ws.Range("E6").Value = 0.469287657961066

$ws.Range("D7").Value = 0.6665986370766313
$ws.Range("E7").Value = 0.3334013629233687

$ws.Range("D8").Value = 0.6996261443360746
$ws.Range("E8").Value = 0.3003738556639254

$ws.Range("D9").Value = 0.5144791887399056
$ws.Range("E9").Value = 0.4855208112600944

$ws.Range("D10").Value = 0.7616316415638997
$ws.Range("E10").Value = 0.2383683584361003

$ws.Range("D11").Value = 0.7251130378056311
$ws.Range("E11").Value = 0.2748869621943689
